# Update the "LOS Galacticos" roster sheet: refresh player/position/team
# data for rows 2-19 (the diff reorders & updates the shared-string backed
# Player / Pozisyon / Takim columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  Player = "Tyus Jones";            Pos = "PG";       Team = "Phoenix Suns" },
    @{ Row = 3;  Player = "James Harden";           Pos = "PG,SG";    Team = "LA Clippers" },
    @{ Row = 4;  Player = "Anthony Edwards";        Pos = "SG,SF";    Team = "Minnesota Timberwolves" },
    @{ Row = 5;  Player = "Fred VanVleet";          Pos = "PG";       Team = "Houston Rockets" },
    @{ Row = 6;  Player = "Amen Thompson";          Pos = "SG,SF";    Team = "Houston Rockets" },
    @{ Row = 7;  Player = "P.J. Washington";        Pos = "PF";       Team = "Dallas Mavericks" },
    @{ Row = 8;  Player = "Bradley Beal";           Pos = "PG,SG,SF"; Team = "Phoenix Suns" },
    @{ Row = 9;  Player = "Zion Williamson";        Pos = "PF,C";     Team = "New Orleans Pelicans" },
    @{ Row = 10; Player = "Ivica Zubac";            Pos = "C";        Team = "LA Clippers" },
    @{ Row = 11; Player = "Jaren Jackson Jr.";      Pos = "PF,C";     Team = "Memphis Grizzlies" },
    @{ Row = 12; Player = "Jayson Tatum";           Pos = "SF,PF";    Team = "Boston Celtics" },
    @{ Row = 13; Player = "Keyonte George";         Pos = "PG,SG";    Team = "Utah Jazz" },
    @{ Row = 14; Player = "Wendell Carter Jr.";     Pos = "PF,C";     Team = "Orlando Magic" },
    @{ Row = 15; Player = "Giannis Antetokounmpo";  Pos = "PF,C";     Team = "Milwaukee Bucks" },
    @{ Row = 16; Player = "Paul George";            Pos = "SG,SF,PF"; Team = "Philadelphia 76ers" },
    @{ Row = 17; Player = "Jonathan Kuminga";       Pos = "SF,PF";    Team = "Golden State Warriors" },
    @{ Row = 18; Player = "Goga Bitadze";           Pos = "C";        Team = "Orlando Magic" },
    @{ Row = 19; Player = "Anfernee Simons";        Pos = "PG,SG";    Team = "Portland Trail Blazers" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value2 = $r.Player
    $ws.Cells.Item($r.Row, 2).Value2 = $r.Pos
    $ws.Cells.Item($r.Row, 3).Value2 = $r.Team
}
